{"js": "// Locate the bullet paragraph that reads exactly \"Trampoline\" (under the\n// \"Mechanics\" heading) and append \" (z_tdm_subway)\" to it, turning it into\n// \"Trampoline (z_tdm_subway)\". This is the only substantive content change\n// in the target revision \u2014 the rest of the diff is just Word's spell-check\n// engine splitting runs / inserting <w:proofErr> markers around words it\n// doesn't recognise (PvE, CoD, gamemode, Bamfuslicator, Planfuslicator,\n// Tangfuslicator, qpath, pullreq'd, ericw, Destructable, yada), which does\n// not change the document's visible text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"Trampoline\") {\n    para.insertText(\" (z_tdm_subway)\", Word.InsertLocation.end);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Locate the bullet paragraph that reads exactly \"Trampoline\" (under the\n# \"Mechanics\" heading) and append \" (z_tdm_subway)\" to it, turning it into\n# \"Trampoline (z_tdm_subway)\". This is the only substantive content change\n# in the target revision -- the rest of the diff is just Word's spell-check\n# engine splitting runs / inserting <w:proofErr> markers around words it\n# doesn't recognise (PvE, CoD, gamemode, Bamfuslicator, Planfuslicator,\n# Tangfuslicator, qpath, pullreq'd, ericw, Destructable, yada), which does\n# not change the document's visible text.\n\n$d = $word.ActiveDocument\n\n$found = $false\nforeach ($p in $d.Paragraphs) {\n  $range = $p.Range\n  $text = $range.Text.Trim()\n  if ($text -eq \"Trampoline\") {\n    $range.MoveEnd(1, -1)       # wdCharacter = 1; trim off the paragraph mark\n    $range.InsertAfter(\" (z_tdm_subway)\")\n    $found = $true\n    break\n  }\n}\n\nif (-not $found) {\n  throw \"Could not find the 'Trampoline' paragraph to update.\"\n}\n"}
